$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the PTP_checkbox / final_k1 / amended_k1 columns (J, K, L) for rows 2-4
# from text "Yes"/"No" values to native boolean TRUE/FALSE values.

$ws.Range("J2").Value = $true
$ws.Range("K2").Value = $true
$ws.Range("L2").Value = $false

$ws.Range("J3").Value = $false
$ws.Range("K3").Value = $false
$ws.Range("L3").Value = $true

$ws.Range("J4").Value = $true
$ws.Range("K4").Value = $true
$ws.Range("L4").Value = $false

# Update the active cell selection to L5
$ws.Range("L5").Select()
